$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 352
$ws1.Range("F4").Value = 4719
$ws1.Range("F6").Value = 476

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 352
$ws4.Range("F4").Value = 4719
$ws4.Range("F8").Value = 476
